$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new "week" column before column B: copy column B and insert the
# copy at column B, which pushes the existing B -> C and C -> D.
$ws.Columns("B:B").Copy()
$ws.Columns("B:B").Insert()

# The freshly inserted column B currently holds a duplicate of the old B
# column (the "Jun_17" header + "UN" rows). Re-label the header for the new
# reporting week; the per-analyst values underneath are already "UN".
$ws.Range("B1").Value = "Jun_26"

# Keep the shifted columns (old "C" now at C/D) at their original 8-wide
# custom width so the layout doesn't change visually.
$ws.Columns("C:C").ColumnWidth = 7.166666666666667
$ws.Columns("D:D").ColumnWidth = 7.166666666666667

# Append the two new tickers/analysts tracked starting this week.
$ws.Range("A28").Value = "Benchmark"
$ws.Range("B28").Value = "UN"

$ws.Range("A29").Value = "Evercore ISI"
$ws.Range("B29").Value = "UN"
